$wb = $excel.ActiveWorkbook

# Duplicate the "Germany" sheet (same layout/styles) and move the copy to the
# very first tab position, then turn it into the new "UK" sheet.
$germany = $wb.Worksheets.Item("Germany")
$germany.Copy($wb.Worksheets.Item(1))
$uk = $wb.Worksheets.Item(1)
$uk.Name = "UK"

# New market-specific content for the UK sheet.
$uk.Range("B2").Value = "UK Market"
$uk.Range("B4").Value = "NGC-3003/T1241/1248/T1257"

# Column widths specific to the UK sheet.
$uk.Columns.Item(1).ColumnWidth = 24.333333333333332
$uk.Columns.Item(2).ColumnWidth = 31
$uk.Columns.Item(3).ColumnWidth = 17.166666666666668
$uk.Columns.Item(4).ColumnWidth = 17.833333333333332

# Match the saved selection on the new sheet.
$null = $uk.Range("B7").Select()
